$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9544904828071594
$ws.Range("B1").Value = 2.028452396392822
$ws.Range("C1").Value = 7.331633567810059
$ws.Range("D1").Value = 2.634469747543335
$ws.Range("E1").Value = 1.403647422790527
